$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 9 de Mayo de 2020 a las 14:04"

# Row 4
$ws.Range("B4").Value = 1322215
$ws.Range("C4").Value = 430
$ws.Range("E4").Value = 1019844
$ws.Range("G4").Value = 7
$ws.Range("H4").Value = 78622

# Row 17
$ws.Range("B17").Value = 59881
$ws.Range("C17").Value = 186
$ws.Range("D17").Value = 17956
$ws.Range("E17").Value = 39935
$ws.Range("G17").Value = 5
$ws.Range("H17").Value = 1990

# Row 25
$ws.Range("B25").Value = 27406
$ws.Range("C25").Value = 138
$ws.Range("D25").Value = 2499
$ws.Range("E25").Value = 23781
$ws.Range("F25").Value = 120
$ws.Range("G25").Value = 12
$ws.Range("H25").Value = 1126

# Row 30
$ws.Range("A30").Value = "Catar"
$ws.Range("B30").Value = 21331
$ws.Range("C30").Value = 1130
$ws.Range("D30").Value = 2449
$ws.Range("E30").Value = 18869
$ws.Range("F30").Value = 72
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 13

# Row 31
$ws.Range("A31").Value = "Bielorrusia"
$ws.Range("B31").Value = 21101
$ws.Range("D31").Value = 5484
$ws.Range("E31").Value = 15496
$ws.Range("F31").Value = 92
$ws.Range("H31").Value = 121

# Row 61
$ws.Range("B61").Value = 4595
$ws.Range("C61").Value = 151
$ws.Range("D61").Value = 2049
$ws.Range("E61").Value = 2538

# Row 84
$ws.Range("A84").Value = "Senegal"
$ws.Range("B84").Value = 1634
$ws.Range("C84").Value = 83
$ws.Range("D84").Value = 643
$ws.Range("E84").Value = 976
$ws.Range("F84").Value = 6
$ws.Range("G84").Value = 2
$ws.Range("H84").Value = 15

# Row 85
$ws.Range("B85").Value = 1622
$ws.Range("C85").Value = 36
$ws.Range("D85").Value = 1112
$ws.Range("E85").Value = 419
$ws.Range("G85").Value = 1
$ws.Range("H85").Value = 91

# Row 86
$ws.Range("A86").Value = "Costa de Marfil"
$ws.Range("B86").Value = 1602
$ws.Range("D86").Value = 754
$ws.Range("E86").Value = 828
$ws.Range("F86").Value = 0
$ws.Range("H86").Value = 20

# Row 108
$ws.Range("A108").Value = "Maldivas"
$ws.Range("B108").Value = 766
$ws.Range("C108").Value = 22
$ws.Range("D108").Value = 20
$ws.Range("E108").Value = 743
$ws.Range("F108").Value = 2
$ws.Range("H108").Value = 3

# Row 109
$ws.Range("A109").Value = "Principado de Andorra"
$ws.Range("B109").Value = 752
$ws.Range("D109").Value = 537
$ws.Range("E109").Value = 168
$ws.Range("F109").Value = 14
$ws.Range("H109").Value = 47

# Row 133
$ws.Range("A133").Value = "Sierra Leona"
$ws.Range("B133").Value = 291
$ws.Range("C133").Value = 34
$ws.Range("D133").Value = 58
$ws.Range("E133").Value = 215
$ws.Range("F133").Value = 0
$ws.Range("G133").Value = 1
$ws.Range("H133").Value = 18

# Row 134
$ws.Range("A134").Value = "Vietnam"
$ws.Range("B134").Value = 288
$ws.Range("D134").Value = 241
$ws.Range("E134").Value = 47
$ws.Range("F134").Value = 8
$ws.Range("H134").Value = 0

# Row 135
$ws.Range("A135").Value = "Congo"
$ws.Range("B135").Value = 274
$ws.Range("D135").Value = 33
$ws.Range("E135").Value = 231
$ws.Range("H135").Value = 10

# Row 136
$ws.Range("A136").Value = "Ruanda"
$ws.Range("B136").Value = 273
$ws.Range("D136").Value = 136
$ws.Range("E136").Value = 137
$ws.Range("H136").Value = 0

# Row 137
$ws.Range("A137").Value = "Republica del Chad"
$ws.Range("B137").Value = 260
$ws.Range("D137").Value = 50
$ws.Range("E137").Value = 182
$ws.Range("H137").Value = 28

# Row 181
$ws.Range("A181").Value = "Zimbabue"
$ws.Range("B181").Value = 35
$ws.Range("C181").Value = 1
$ws.Range("D181").Value = 9
$ws.Range("E181").Value = 22
$ws.Range("H181").Value = 4

# Row 182
$ws.Range("A182").Value = "Yemen"
$ws.Range("D182").Value = 1
$ws.Range("E182").Value = 26
$ws.Range("H182").Value = 7

# Row 212
$ws.Range("A212").Value = "Islas Virgenes Britanicas"
$ws.Range("D212").Value = 4
$ws.Range("H212").Value = 1

# Row 213
$ws.Range("A213").Value = "Butan"
$ws.Range("D213").Value = 5
$ws.Range("H213").Value = 0
